# Release Aspose.Cells Cloud SDK 23.12 - refresh TestData/NewCopy.xlsx
# Sheet1 gets a value, two new (empty) worksheets are appended.

$wb = $excel.ActiveWorkbook

# --- Sheet1: enter a value and leave the cursor on the next cell down ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A1").Value = 1111
$ws1.PageSetup.PaperSize = 9
$ws1.PageSetup.Orientation = 1
$ws1.Range("A2").Select() | Out-Null

# --- Append Sheet2 and Sheet3 after Sheet1, in order ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"
$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1

$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "Sheet3"
$ws3.PageSetup.PaperSize = 9
$ws3.PageSetup.Orientation = 1

# --- keep Sheet1 as the active/selected sheet ---
$ws1.Select() | Out-Null
